$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A values are numeric-looking Cow IDs that must be retained as
# zero-padded text (e.g. "07596"), so force text format before writing.
$ws.Range("A2:A51").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "15677"
$ws.Range("B2").Value = 0

# Row 3
$ws.Range("A3").Value = "60864"

# Row 4
$ws.Range("A4").Value = "81875"

# Row 5
$ws.Range("A5").Value = "23389"

# Row 6
$ws.Range("A6").Value = "13829"

# Row 7
$ws.Range("A7").Value = "07596"
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "Not being milked due to clinical mastitis"

# Row 8
$ws.Range("A8").Value = "40285"

# Row 9
$ws.Range("A9").Value = "28571"

# Row 10
$ws.Range("A10").Value = "05305"

# Row 11
$ws.Range("A11").Value = "20656"
$ws.Range("B11").Value = 0

# Row 12
$ws.Range("A12").Value = "35943"
$ws.Range("B12").Value = 1

# Row 13
$ws.Range("A13").Value = "78505"

# Row 14
$ws.Range("A14").Value = "40654"

# Row 15
$ws.Range("A15").Value = "70233"

# Row 16
$ws.Range("A16").Value = "26654"

# Row 17
$ws.Range("A17").Value = "80194"

# Row 18
$ws.Range("A18").Value = "80614"
$ws.Range("B18").Value = 1

# Row 19
$ws.Range("A19").Value = "30708"
$ws.Range("B19").Value = 0
$ws.Range("C19").ClearContents()

# Row 20
$ws.Range("A20").Value = "13011"

# Row 21
$ws.Range("A21").Value = "81448"
$ws.Range("B21").Value = 0

# Row 22
$ws.Range("A22").Value = "30324"

# Row 23
$ws.Range("A23").Value = "93218"
$ws.Range("B23").ClearContents()
$ws.Range("C23").Value = "Not being milked due to clinical mastitis"

# Row 24
$ws.Range("A24").Value = "58817"

# Row 25
$ws.Range("A25").Value = "21603"

# Row 26
$ws.Range("A26").Value = "84045"

# Row 27
$ws.Range("A27").Value = "98492"
$ws.Range("B27").Value = 0

# Row 28
$ws.Range("A28").Value = "14116"

# Row 29
$ws.Range("A29").Value = "67775"

# Row 30
$ws.Range("A30").Value = "16065"

# Row 31
$ws.Range("A31").Value = "51579"

# Row 32
$ws.Range("A32").Value = "00194"

# Row 33
$ws.Range("A33").Value = "96132"
$ws.Range("B33").ClearContents()
$ws.Range("C33").Value = "Walked in front of another cow"

# Row 34
$ws.Range("A34").Value = "50774"

# Row 35
$ws.Range("A35").Value = "99302"

# Row 36
$ws.Range("A36").Value = "50895"

# Row 37
$ws.Range("A37").Value = "42712"

# Row 38
$ws.Range("A38").Value = "04793"

# Row 39
$ws.Range("A39").Value = "96195"

# Row 40
$ws.Range("A40").Value = "93604"

# Row 41
$ws.Range("A41").Value = "27887"

# Row 42
$ws.Range("A42").Value = "57635"

# Row 43
$ws.Range("A43").Value = "51078"

# Row 44
$ws.Range("A44").Value = "45360"
$ws.Range("B44").Value = 0
$ws.Range("C44").ClearContents()

# Row 45
$ws.Range("A45").Value = "71880"

# Row 46
$ws.Range("A46").Value = "34504"

# Row 47
$ws.Range("A47").Value = "23926"

# Row 48
$ws.Range("A48").Value = "09833"

# Row 49
$ws.Range("A49").Value = "43064"

# Row 50
$ws.Range("A50").Value = "25163"

# Row 51
$ws.Range("A51").Value = "33050"

